$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: same header style as B1:E1 (style index 1 -> bold/border/centered)
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

# time_taken values for rows 2..244 (data rows), no special style
$timeTaken = @(
    "2021-10-05 13:40:08.829838",
    "2021-10-05 13:40:08.829849",
    "2021-10-05 13:40:08.829852",
    "2021-10-05 13:40:08.829855",
    "2021-10-05 13:40:08.829858",
    "2021-10-05 13:40:08.829861",
    "2021-10-05 13:40:08.829863",
    "2021-10-05 13:40:08.829866",
    "2021-10-05 13:40:08.829868",
    "2021-10-05 13:40:08.829871",
    "2021-10-05 13:40:08.829873",
    "2021-10-05 13:40:08.829876",
    "2021-10-05 13:40:08.829878",
    "2021-10-05 13:40:08.829881",
    "2021-10-05 13:40:08.829883",
    "2021-10-05 13:40:08.829886",
    "2021-10-05 13:40:08.829889",
    "2021-10-05 13:40:08.829891",
    "2021-10-05 13:40:08.829894",
    "2021-10-05 13:40:08.829896",
    "2021-10-05 13:40:08.829899",
    "2021-10-05 13:40:08.829901",
    "2021-10-05 13:40:08.829904",
    "2021-10-05 13:40:08.829906",
    "2021-10-05 13:40:08.829909",
    "2021-10-05 13:40:08.829911",
    "2021-10-05 13:40:08.829914",
    "2021-10-05 13:40:08.829916",
    "2021-10-05 13:40:08.829919",
    "2021-10-05 13:40:08.829922",
    "2021-10-05 13:40:08.829924",
    "2021-10-05 13:40:08.829927",
    "2021-10-05 13:40:08.829929",
    "2021-10-05 13:40:08.829932",
    "2021-10-05 13:40:08.829935",
    "2021-10-05 13:40:08.829937",
    "2021-10-05 13:40:08.829939",
    "2021-10-05 13:40:08.829942",
    "2021-10-05 13:40:08.829944",
    "2021-10-05 13:40:08.829947",
    "2021-10-05 13:40:08.829950",
    "2021-10-05 13:40:08.829952",
    "2021-10-05 13:40:08.829955",
    "2021-10-05 13:40:08.829957",
    "2021-10-05 13:40:08.829960",
    "2021-10-05 13:40:08.829962",
    "2021-10-05 13:40:08.829965",
    "2021-10-05 13:40:08.829967",
    "2021-10-05 13:40:08.829970",
    "2021-10-05 13:40:08.829972",
    "2021-10-05 13:40:08.829975",
    "2021-10-05 13:40:08.829977",
    "2021-10-05 13:40:08.829980",
    "2021-10-05 13:40:08.829983",
    "2021-10-05 13:40:08.829985",
    "2021-10-05 13:40:08.829988",
    "2021-10-05 13:40:08.829990",
    "2021-10-05 13:40:08.829993",
    "2021-10-05 13:40:08.829995",
    "2021-10-05 13:40:08.829998",
    "2021-10-05 13:40:08.830000",
    "2021-10-05 13:40:08.830003",
    "2021-10-05 13:40:08.830005",
    "2021-10-05 13:40:08.830008",
    "2021-10-05 13:40:08.830011",
    "2021-10-05 13:40:08.830014",
    "2021-10-05 13:40:08.830017",
    "2021-10-05 13:40:08.830019",
    "2021-10-05 13:40:08.830022",
    "2021-10-05 13:40:08.830024",
    "2021-10-05 13:40:08.830027",
    "2021-10-05 13:40:08.830029",
    "2021-10-05 13:40:08.830032",
    "2021-10-05 13:40:08.830034",
    "2021-10-05 13:40:08.830036",
    "2021-10-05 13:40:08.830039",
    "2021-10-05 13:40:08.830043",
    "2021-10-05 13:40:08.830046",
    "2021-10-05 13:40:08.830049",
    "2021-10-05 13:40:08.830052",
    "2021-10-05 13:40:08.830054",
    "2021-10-05 13:40:08.830057",
    "2021-10-05 13:40:08.830059",
    "2021-10-05 13:40:08.830062",
    "2021-10-05 13:40:08.830065",
    "2021-10-05 13:40:08.830067",
    "2021-10-05 13:40:08.830070",
    "2021-10-05 13:40:08.830072",
    "2021-10-05 13:40:08.830075",
    "2021-10-05 13:40:08.830078",
    "2021-10-05 13:40:08.830080",
    "2021-10-05 13:40:08.830083",
    "2021-10-05 13:40:08.830086",
    "2021-10-05 13:40:08.830089",
    "2021-10-05 13:40:08.830092",
    "2021-10-05 13:40:08.830094",
    "2021-10-05 13:40:08.830097",
    "2021-10-05 13:40:08.830099",
    "2021-10-05 13:40:08.830102",
    "2021-10-05 13:40:08.830105",
    "2021-10-05 13:40:08.830107",
    "2021-10-05 13:40:08.830110",
    "2021-10-05 13:40:08.830112",
    "2021-10-05 13:40:08.830115",
    "2021-10-05 13:40:08.830117",
    "2021-10-05 13:40:08.830120",
    "2021-10-05 13:40:08.830122",
    "2021-10-05 13:40:08.830125",
    "2021-10-05 13:40:08.830129",
    "2021-10-05 13:40:08.830132",
    "2021-10-05 13:40:08.830135",
    "2021-10-05 13:40:08.830137",
    "2021-10-05 13:40:08.830140",
    "2021-10-05 13:40:08.830142",
    "2021-10-05 13:40:08.830145",
    "2021-10-05 13:40:08.830147",
    "2021-10-05 13:40:08.830149",
    "2021-10-05 13:40:08.830152",
    "2021-10-05 13:40:08.830154",
    "2021-10-05 13:40:08.830157",
    "2021-10-05 13:40:08.830159",
    "2021-10-05 13:40:08.830162",
    "2021-10-05 13:40:08.830164",
    "2021-10-05 13:40:08.830167",
    "2021-10-05 13:40:08.830169",
    "2021-10-05 13:40:08.830172",
    "2021-10-05 13:40:08.830174",
    "2021-10-05 13:40:08.830177",
    "2021-10-05 13:40:08.830181",
    "2021-10-05 13:40:08.830184",
    "2021-10-05 13:40:08.830187",
    "2021-10-05 13:40:08.830189",
    "2021-10-05 13:40:08.830192",
    "2021-10-05 13:40:08.830195",
    "2021-10-05 13:40:08.830197",
    "2021-10-05 13:40:08.830200",
    "2021-10-05 13:40:08.830202",
    "2021-10-05 13:40:08.830205",
    "2021-10-05 13:40:08.830207",
    "2021-10-05 13:40:08.830210",
    "2021-10-05 13:40:08.830212",
    "2021-10-05 13:40:08.830215",
    "2021-10-05 13:40:08.830217",
    "2021-10-05 13:40:08.830219",
    "2021-10-05 13:40:08.830222",
    "2021-10-05 13:40:08.830225",
    "2021-10-05 13:40:08.830227",
    "2021-10-05 13:40:08.830230",
    "2021-10-05 13:40:08.830232",
    "2021-10-05 13:40:08.830235",
    "2021-10-05 13:40:08.830237",
    "2021-10-05 13:40:08.830240",
    "2021-10-05 13:40:08.830242",
    "2021-10-05 13:40:08.830245",
    "2021-10-05 13:40:08.830247",
    "2021-10-05 13:40:08.830250",
    "2021-10-05 13:40:08.830252",
    "2021-10-05 13:40:08.830255",
    "2021-10-05 13:40:08.830257",
    "2021-10-05 13:40:08.830260",
    "2021-10-05 13:40:08.830262",
    "2021-10-05 13:40:08.830265",
    "2021-10-05 13:40:08.830267",
    "2021-10-05 13:40:08.830270",
    "2021-10-05 13:40:08.830272",
    "2021-10-05 13:40:08.830275",
    "2021-10-05 13:40:08.830277",
    "2021-10-05 13:40:08.830280",
    "2021-10-05 13:40:08.830282",
    "2021-10-05 13:40:08.830285",
    "2021-10-05 13:40:08.830287",
    "2021-10-05 13:40:08.830290",
    "2021-10-05 13:40:08.830294",
    "2021-10-05 13:40:08.830297",
    "2021-10-05 13:40:08.830299",
    "2021-10-05 13:40:08.830302",
    "2021-10-05 13:40:08.830304",
    "2021-10-05 13:40:08.830307",
    "2021-10-05 13:40:08.830309",
    "2021-10-05 13:40:08.830312",
    "2021-10-05 13:40:08.830314",
    "2021-10-05 13:40:08.830317",
    "2021-10-05 13:40:08.830319",
    "2021-10-05 13:40:08.830322",
    "2021-10-05 13:40:08.830324",
    "2021-10-05 13:40:08.830327",
    "2021-10-05 13:40:08.830329",
    "2021-10-05 13:40:08.830332",
    "2021-10-05 13:40:08.830335",
    "2021-10-05 13:40:08.830337",
    "2021-10-05 13:40:08.830340",
    "2021-10-05 13:40:08.830342",
    "2021-10-05 13:40:08.830344",
    "2021-10-05 13:40:08.830347",
    "2021-10-05 13:40:08.830350",
    "2021-10-05 13:40:08.830352",
    "2021-10-05 13:40:08.830355",
    "2021-10-05 13:40:08.830357",
    "2021-10-05 13:40:08.830360",
    "2021-10-05 13:40:08.830362",
    "2021-10-05 13:40:08.830365",
    "2021-10-05 13:40:08.830368",
    "2021-10-05 13:40:08.830370",
    "2021-10-05 13:40:08.830373",
    "2021-10-05 13:40:08.830375",
    "2021-10-05 13:40:08.830378",
    "2021-10-05 13:40:08.830380",
    "2021-10-05 13:40:08.830383",
    "2021-10-05 13:40:08.830385",
    "2021-10-05 13:40:08.830388",
    "2021-10-05 13:40:08.830390",
    "2021-10-05 13:40:08.830393",
    "2021-10-05 13:40:08.830395",
    "2021-10-05 13:40:08.830398",
    "2021-10-05 13:40:08.830401",
    "2021-10-05 13:40:08.830403",
    "2021-10-05 13:40:08.830406",
    "2021-10-05 13:40:08.830408",
    "2021-10-05 13:40:08.830410",
    "2021-10-05 13:40:08.830413",
    "2021-10-05 13:40:08.830415",
    "2021-10-05 13:40:08.830418",
    "2021-10-05 13:40:08.830420",
    "2021-10-05 13:40:08.830423",
    "2021-10-05 13:40:08.830425",
    "2021-10-05 13:40:08.830428",
    "2021-10-05 13:40:08.830430",
    "2021-10-05 13:40:08.830433",
    "2021-10-05 13:40:08.830435",
    "2021-10-05 13:40:08.830438",
    "2021-10-05 13:40:08.830440",
    "2021-10-05 13:40:08.830443",
    "2021-10-05 13:40:08.830447",
    "2021-10-05 13:40:08.830450",
    "2021-10-05 13:40:08.830453",
    "2021-10-05 13:40:08.830455",
    "2021-10-05 13:40:08.830458",
    "2021-10-05 13:40:08.830461",
    "2021-10-05 13:40:08.830463",
    "2021-10-05 13:40:08.830466",
    "2021-10-05 13:40:08.830469",
    "2021-10-05 13:40:08.830471",
    "2021-10-05 13:40:08.830474"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Host "done"